$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "Datos actualizados" timestamp label (A1) ---
$ws.Range("A1").Value2 = "Datos actualizados a 18 de Julio de 2020 a las 15:18"

# --- Update country statistics rows ---

# Row 6: India
$ws.Range("B6").Value2 = 1047238
$ws.Range("C6").Value2 = 6781
$ws.Range("D6").Value2 = 656464
$ws.Range("E6").Value2 = 364421
$ws.Range("G6").Value2 = 68
$ws.Range("H6").Value2 = 26353

# Row 16: Arabia Saudita
$ws.Range("B16").Value2 = 248416
$ws.Range("C16").Value2 = 2565
$ws.Range("D16").Value2 = 194218
$ws.Range("E16").Value2 = 51751
$ws.Range("G16").Value2 = 40
$ws.Range("H16").Value2 = 2447

# Row 23: Argentina
$ws.Range("D23").Value2 = 52607
$ws.Range("E23").Value2 = 64490
$ws.Range("G23").Value2 = 26
$ws.Range("H23").Value2 = 2204

# Row 26: Irak
$ws.Range("B26").Value2 = 90220
$ws.Range("C26").Value2 = 2049
$ws.Range("D26").Value2 = 58492
$ws.Range("E26").Value2 = 28037
$ws.Range("G26").Value2 = 75
$ws.Range("H26").Value2 = 3691

# Rows 61/62: Moldavia & Serbia swap position (Serbia moves above Moldavia)
# and Serbia gets new updated figures while Moldavia keeps its old figures.
$ws.Range("A61").Value2 = "Serbia"
$ws.Range("B61").Value2 = 20498
$ws.Range("C61").Value2 = 389
$ws.Range("D61").Value2 = 14047
$ws.Range("E61").Value2 = 5990
$ws.Range("G61").Value2 = 9
$ws.Range("H61").Value2 = 461

$ws.Range("A62").Value2 = "Moldavia"
$ws.Range("B62").Value2 = 20494
$ws.Range("C62").Value2 = 0
$ws.Range("D62").Value2 = 13913
$ws.Range("E62").Value2 = 5906
$ws.Range("G62").Value2 = 0
$ws.Range("H62").Value2 = 675

# Row 100: Croacia
$ws.Range("B100").Value2 = 4253
$ws.Range("C100").Value2 = 116
$ws.Range("D100").Value2 = 2929
$ws.Range("E100").Value2 = 1204

# Row 188: Liechtenstein
$ws.Range("B188").Value2 = 86
$ws.Range("C188").Value2 = 1
$ws.Range("E188").Value2 = 4

$wb.Save()
